$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("priority") before the existing jira_id column,
# shifting jira_id from column D to column E.
$ws.Columns.Item(4).Insert()

# Header for the new column
$ws.Range("D1").Value = "priority"

# Move the jira_id value for row 2 (previously in D2) into the new E2 cell,
# and set the priority value for that row.
$ws.Range("E2").Value = "TSET-1"
$ws.Range("D2").Value = "high"
